$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.311.30'
$ws.Range('E2').Value = '  -1.18%  '
$ws.Range('D3').Value = '1.588.75'
$ws.Range('E3').Value = '  -0.52%  '
$ws.Range('E4').Value = '  -0.37%  '
$ws.Range('D5').Value = '''209.97'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.28%  '
$ws.Range('D6').Value = '''0.504'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.98%  '
$ws.Range('E7').Value = '  -0.34%  '
$ws.Range('E8').Value = '  -0.77%  '
$ws.Range('E9').Value = '  -0.44%  '
$ws.Range('E10').Value = '  -0.45%  '
$ws.Range('D12').Value = '1.810.88'
$ws.Range('E12').Value = '  -0.58%  '
$ws.Range('D13').Value = '''4.07'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.64%  '
$ws.Range('D14').Value = '1.560.94'
$ws.Range('E14').Value = '  -2.25%  '
$ws.Range('D16').Value = '''64.33'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.39%  '
$ws.Range('D17').Value = '26.319.28'
$ws.Range('D18').Value = '0.0₃0727'
$ws.Range('E18').Value = '  -1.62%  '
$ws.Range('D19').Value = '''7.48'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +5.65%  '
$ws.Range('D20').Value = '''211.03'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.34%  '
$ws.Range('E21').Value = '  -0.31%  '
$ws.Range('E22').Value = '  -0.35%  '
$ws.Range('E23').Value = '  -2.79%  '
$ws.Range('D24').Value = '''8.95'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('D25').Value = '''144.93'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.82%  '
$ws.Range('E26').Value = '  -0.33%  '
$ws.Range('E28').Value = '  -0.43%  '
$ws.Range('E29').Value = '  -0.15%  '
$ws.Range('E30').Value = '  +0.07%  '
$ws.Range('E31').Value = '  -0.15%  '
$ws.Range('E32').Value = '  -1.16%  '
$ws.Range('E33').Value = '  +0.90%  '
$ws.Range('D34').Value = '1.305.22'
$ws.Range('E34').Value = '  +2.13%  '
$ws.Range('E35').Value = '  +2.00%  '
$ws.Range('D36').Value = '''2.44'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -1.64%  '
$ws.Range('E38').Value = '  +0.25%  '
$ws.Range('E39').Value = '  -11.54%  '
$ws.Range('E40').Value = '  -1.96%  '
$ws.Range('E41').Value = '  -0.34%  '
$ws.Range('E42').Value = '  +3.43%  '
$ws.Range('E43').Value = '  -0.52%  '
$ws.Range('E44').Value = '  -1.34%  '
$ws.Range('D45').Value = '''62.37'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.45%  '
$ws.Range('D46').Value = '1.723.81'
$ws.Range('E46').Value = '  -0.51%  '
$ws.Range('D47').Value = '''87.73'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.92%  '
$ws.Range('D48').Value = '''1.49'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -5.35%  '
$ws.Range('E49').Value = '  -1.45%  '
$ws.Range('E50').Value = '  -4.48%  '
$ws.Range('E51').Value = '  -0.36%  '
